# Weekly update: insert a new "Acelga" price record at row 441 for
# Macroferia Regional de Talca. All existing records at row 441 and below
# shift down by one row (this mirrors how the source workbook appends the
# latest weekly reading at the top of the date-ordered block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 441; rows 441:516 become 442:517 and pick up the
# row's existing formatting (e.g. the date style on column D).
$ws.Rows.Item(441).Insert()

# Populate the newly inserted row with this week's reading.
$ws.Range("A441").Value = 5
$ws.Range("B441").Value = "Macroferia Regional de Talca"
$ws.Range("C441").Value = "Maule"
$ws.Range("D441").Value = 45218
$ws.Range("E441").Value = 7
$ws.Range("F441").Value = 100112009
$ws.Range("G441").Value = "Acelga"
$ws.Range("H441").Value = "Sin especificar"
$ws.Range("I441").Value = "Primera"
$ws.Range("J441").Value = 400
$ws.Range("K441").Value = 2000
$ws.Range("L441").Value = 2000
$ws.Range("M441").Value = 2000
$ws.Range("N441").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O441").Value = "Región del Maule"
$ws.Range("P441").Value = 500
$ws.Range("Q441").Value = 4
$ws.Range("R441").Value = "Hortaliza"
